$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '28.019.89'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '1.892.22'
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +1.34%  '

$ws.Range("D5").Value = '336.24'

$ws.Range("E6").Value = '  +1.26%  '

$ws.Range("D7").Value = '0.4713'
$ws.Range("E7").Value = '  -0.43%  '

$ws.Range("D8").Value = '0.3954'
$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("D9").Value = '46.84'
$ws.Range("E9").Value = '  -2.87%  '

$ws.Range("D10").Value = '0.08025'
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("D12").Value = '21.83'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").Value = '1.888.83'
$ws.Range("E13").Value = '  -1.40%  '

$ws.Range("E14").Value = '  +0.55%  '

$ws.Range("D15").Value = '7.181'
$ws.Range("E15").Value = '  -0.29%  '

$ws.Range("D16").Value = '1.018'
$ws.Range("E16").Value = '  +1.47%  '

$ws.Range("D17").Value = '0.06780'
$ws.Range("E17").Value = '  +2.35%  '

$ws.Range("D18").Value = '88.05'
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").Value = '0.00001053'
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("D20").Value = '17.22'
$ws.Range("E20").Value = '  -0.93%  '

$ws.Range("E21").Value = '  +1.30%  '

$ws.Range("D22").Value = '28.021.59'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = '5.509'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("D25").Value = '2.364'
$ws.Range("E25").Value = '  +2.07%  '

$ws.Range("D26").Value = '2.115.67'
$ws.Range("E26").Value = '  -1.10%  '

$ws.Range("D27").Value = '159.52'
$ws.Range("E27").Value = '  +1.16%  '

$ws.Range("D28").Value = '20.05'
$ws.Range("E28").Value = '  -1.12%  '

$ws.Range("D29").Value = '2.106'
$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("D30").Value = '5.502'
$ws.Range("E30").Value = '  -2.17%  '

$ws.Range("D31").Value = '121.67'
$ws.Range("E31").Value = '  -0.81%  '

$ws.Range("D32").Value = '0.09575'
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").Value = '0.9670'
$ws.Range("E33").Value = '  -1.78%  '

$ws.Range("E34").Value = '  +0.84%  '

$ws.Range("D35").Value = '5.363'
$ws.Range("E35").Value = '  +0.62%  '

$ws.Range("D36").Value = '1.365'
$ws.Range("E36").Value = '  -6.78%  '

$ws.Range("D37").Value = '0.06134'
$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").Value = '0.02253'
$ws.Range("E38").Value = '  -0.53%  '

$ws.Range("D39").Value = '1.216'
$ws.Range("E39").Value = '  -1.23%  '

$ws.Range("D40").Value = '8.238'
$ws.Range("E40").Value = '  -0.18%  '

$ws.Range("D41").Value = '0.5979'
$ws.Range("E41").Value = '  -0.90%  '

$ws.Range("D42").Value = '0.1906'
$ws.Range("E42").Value = '  +0.09%  '

$ws.Range("D43").Value = '10.35'
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").Value = '1.266'
$ws.Range("E44").Value = '  -0.56%  '

$ws.Range("D45").Value = '0.5711'
$ws.Range("E45").Value = '  -0.18%  '

$ws.Range("D46").Value = '12.19'
$ws.Range("E46").Value = '  -1.05%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '1.948'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.396'
$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("D49").Value = '0.06869'
$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("D50").Value = '113.72'
$ws.Range("E50").Value = '  -0.09%  '

$ws.Range("E51").Value = '  -0.49%  '

